$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 264; this shifts all existing rows 264..376
# down to 265..377, matching the target diff (dimension grows to A1:R377).
$ws.Rows(264).Insert()

# Populate the newly inserted row 264 with the new price-record data.
$ws.Range("A264").Value = 6
$ws.Range("B264").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C264").Value = "Metropolitana"
$ws.Range("D264").Value = 44510
$ws.Range("E264").Value = 13
$ws.Range("F264").Value = 100112044
$ws.Range("G264").Value = "Perejil"
$ws.Range("H264").Value = "Sin especificar"
$ws.Range("I264").Value = "Primera"
$ws.Range("J264").Value = 180
$ws.Range("K264").Value = 12000
$ws.Range("L264").Value = 13000
$ws.Range("M264").Value = 12611
$ws.Range("N264").Value = "$/docena de atados"
$ws.Range("O264").Value = "Región Metropolitana"
$ws.Range("P264").Value = 4204
$ws.Range("Q264").Value = 3
$ws.Range("R264").Value = "Hortaliza"
